$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-11-15 16:40:57"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-11-15 16:40:42"
$wsZhCn.Range("K4").Value = "2016-11-15 16:41:43"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-11-15 16:40:57"
$wsDeDe.Range("K4").Value = "2016-11-15 16:42:01"
